# The docx has a "first page" header/footer pair (header1.xml/footer1.xml)
# and a "default" header/footer pair (header2.xml/footer2.xml) used on the
# remaining pages. In the Word object model, HeadersFooters.Item(1) is the
# primary/default story and Item(2) is the first-page story, so:
#   Headers.Item(1) -> header2.xml   Headers.Item(2) -> header1.xml
#   Footers.Item(1) -> footer2.xml   Footers.Item(2) -> footer1.xml
# Each story holds a single inline picture whose <wp:docPr>/<pic:cNvPr>
# "name" attribute needs to be swapped (the ids/descr stay untouched).
# InlineShape has no settable Name property, so the picture's drawing XML
# is rebuilt (same geometry/ids/descr/relationship) via Range.InsertXML,
# which is the supported way to rewrite a range's underlying OOXML.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

$pkgNs = 'xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"'
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$wpNs = 'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"'
$aNs = 'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"'
$picNs = 'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"'
$rNs = 'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

function Rename-InlinePicture($range, $descr, $picId, $cx, $cy, $newName) {
    $drawing = '<w:drawing ' + $wNs + ' ' + $wpNs + ' ' + $aNs + ' ' + $picNs + ' ' + $rNs + '>' +
        '<wp:inline distB="0" distT="0" distL="0" distR="0">' +
            '<wp:extent cx="' + $cx + '" cy="' + $cy + '"/>' +
            '<wp:effectExtent b="0" l="0" r="0" t="0"/>' +
            '<wp:docPr descr="' + $descr + '" id="' + $picId + '" name="' + $newName + '"/>' +
            '<a:graphic>' +
                '<a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
                    '<pic:pic>' +
                        '<pic:nvPicPr>' +
                            '<pic:cNvPr descr="' + $descr + '" id="0" name="' + $newName + '"/>' +
                            '<pic:cNvPicPr preferRelativeResize="0"/>' +
                        '</pic:nvPicPr>' +
                        '<pic:blipFill>' +
                            '<a:blip r:embed="rId1"/>' +
                            '<a:srcRect b="0" l="0" r="0" t="0"/>' +
                            '<a:stretch><a:fillRect/></a:stretch>' +
                        '</pic:blipFill>' +
                        '<pic:spPr>' +
                            '<a:xfrm><a:off x="0" y="0"/><a:ext cx="' + $cx + '" cy="' + $cy + '"/></a:xfrm>' +
                            '<a:prstGeom prst="rect"/>' +
                            '<a:ln/>' +
                        '</pic:spPr>' +
                    '</pic:pic>' +
                '</a:graphicData>' +
            '</a:graphic>' +
        '</wp:inline>' +
    '</w:drawing>'

    $xmlFrag = '<pkg:package ' + $pkgNs + '>' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
                '<w:document ' + $wNs + ' ' + $wpNs + ' ' + $aNs + ' ' + $picNs + ' ' + $rNs + '>' +
                    '<w:body><w:p><w:r>' + $drawing + '</w:r></w:p></w:body>' +
                '</w:document>' +
            '</pkg:xmlData>' +
        '</pkg:part>' +
    '</pkg:package>'

    $range.InsertXML($xmlFrag)
}

# --- Headers: BTec_Logo-Orange JPEG, image2.jpg -> image1.jpg ---
$hdrDefault = $sec.Headers.Item(1)   # header2.xml, docPr id="3"
$shp = $hdrDefault.Range.InlineShapes.Item(1)
Rename-InlinePicture $shp.Range "BTec_Logo-Orange" 3 914400 277792 "image1.jpg"

$hdrFirst = $sec.Headers.Item(2)     # header1.xml, docPr id="1"
$shp = $hdrFirst.Range.InlineShapes.Item(1)
Rename-InlinePicture $shp.Range "BTec_Logo-Orange" 1 914400 277792 "image1.jpg"

# --- Footers: Pearson logo PNG, image1.png -> image2.png ---
$pearsonDescr = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"

$ftrDefault = $sec.Footers.Item(1)   # footer2.xml, docPr id="4"
$shp = $ftrDefault.Range.InlineShapes.Item(1)
Rename-InlinePicture $shp.Range $pearsonDescr 4 952500 285750 "image2.png"

$ftrFirst = $sec.Footers.Item(2)     # footer1.xml, docPr id="2"
$shp = $ftrFirst.Range.InlineShapes.Item(1)
Rename-InlinePicture $shp.Range $pearsonDescr 2 952500 285750 "image2.png"

Write-Output "done"
